$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price (D) and Volume (E) columns so numeric-looking
# strings like "233.75" or "21.05" are not reinterpreted as numbers, and
# percentage strings keep their surrounding whitespace.
$priceRange = $ws.Range("D2:D51")
$volRange = $ws.Range("E2:E51")
$priceRange.NumberFormat = "@"
$volRange.NumberFormat = "@"

$ws.Range('D2').Value = '37.467.32'
$ws.Range('E2').Value = '  +2.23%  '
$ws.Range('D3').Value = '2.068.65'
$ws.Range('E3').Value = '  +2.33%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '233.75'
$ws.Range('E5').Value = '  -0.53%  '
$ws.Range('D6').Value = '0.620'
$ws.Range('E6').Value = '  +3.18%  '
$ws.Range('D7').Value = '58.13'
$ws.Range('E7').Value = '  +5.83%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '0.385'
$ws.Range('E9').Value = '  +3.90%  '
$ws.Range('D10').Value = '59.08'
$ws.Range('E10').Value = '  +1.88%  '
$ws.Range('D11').Value = '0.0763'
$ws.Range('E11').Value = '  +1.80%  '
$ws.Range('E12').Value = '  -0.14%  '
$ws.Range('D13').Value = '2.375.20'
$ws.Range('E13').Value = '  +2.19%  '
$ws.Range('D14').Value = '14.43'
$ws.Range('E14').Value = '  +1.88%  '
$ws.Range('D15').Value = '21.05'
$ws.Range('E15').Value = '  +4.32%  '
$ws.Range('D16').Value = '0.776'
$ws.Range('E16').Value = '  +1.77%  '
$ws.Range('E17').Value = '  +2.04%  '
$ws.Range('D18').Value = '2.066.14'
$ws.Range('E18').Value = '  +2.10%  '
$ws.Range('D19').Value = '37.656.90'
$ws.Range('E19').Value = '  +1.19%  '
$ws.Range('D20').Value = '6.17'
$ws.Range('E20').Value = '  +14.80%  '
$ws.Range('D21').Value = '69.36'
$ws.Range('E21').Value = '  +2.32%  '
$ws.Range('D22').Value = '0.0₃0813'
$ws.Range('E22').Value = '  +1.96%  '
$ws.Range('D23').Value = '226.15'
$ws.Range('E23').Value = '  +2.50%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').Value = '2.45'
$ws.Range('E25').Value = '  +2.12%  '
$ws.Range('E26').Value = '  +0.72%  '
$ws.Range('D27').Value = '165.99'
$ws.Range('E27').Value = '  +1.45%  '
$ws.Range('D28').Value = '1.49'
$ws.Range('E28').Value = '  +7.47%  '
$ws.Range('D29').Value = '8.96'
$ws.Range('E29').Value = '  +3.68%  '
$ws.Range('E30').Value = '  +1.47%  '
$ws.Range('D31').Value = '19.19'
$ws.Range('E31').Value = '  +1.84%  '
$ws.Range('E32').Value = '  +0.81%  '
$ws.Range('D33').Value = '4.49'
$ws.Range('E33').Value = '  +2.84%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = '2.59'
$ws.Range('E34').Value = '  +5.96%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.0624'
$ws.Range('E35').Value = '  +3.30%  '
$ws.Range('D36').Value = '4.60'
$ws.Range('E36').Value = '  +8.17%  '
$ws.Range('D37').Value = '6.02'
$ws.Range('E37').Value = '  +4.77%  '
$ws.Range('E38').Value = '  -0.20%  '
$ws.Range('D39').Value = '3.32'
$ws.Range('E39').Value = '  +0.80%  '
$ws.Range('E40').Value = '  -0.79%  '
$ws.Range('D41').Value = '4.70'
$ws.Range('E41').Value = '  +14.47%  '
$ws.Range('E42').Value = '  -0.87%  '
$ws.Range('E43').Value = '  +2.35%  '
$ws.Range('D44').Value = '1.459.51'
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('E45').Value = '  +6.12%  '
$ws.Range('D46').Value = '95.56'
$ws.Range('E46').Value = '  +6.00%  '
$ws.Range('E47').Value = '  +3.84%  '
$ws.Range('D48').Value = '15.72'
$ws.Range('E48').Value = '  +2.75%  '
$ws.Range('E49').Value = '  +2.43%  '
$ws.Range('D50').Value = '7.19'
$ws.Range('E50').Value = '  +4.77%  '
$ws.Range('E51').Value = '  +2.16%  '

# Restore the default (unstyled) cell style so we do not introduce new
# formatting that was not part of the original workbook.
$priceRange.Style = "Normal"
$volRange.Style = "Normal"
